# Added ifo GDP component analysis preprocessing: refreshed the naive
# QoQ error series (GVA matched to ifoCAST since 2021) with updated
# averaged error values, including newly populated trailing-diagonal
# cells for the later quarters (rows 11-20).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -0.8230868241805323
$ws.Range("G2").Value = -0.10215989263453
$ws.Range("H2").Value = -0.8111118172483847
$ws.Range("I2").Value = -0.3190058249259076
$ws.Range("J2").Value = 0.1507989223470795
$ws.Range("K2").Value = 0.4636241173125255

$ws.Range("E3").Value = -0.931701356396673
$ws.Range("F3").Value = -0.2107744248506708
$ws.Range("G3").Value = -0.9197263494645256
$ws.Range("H3").Value = -0.4276203571420484
$ws.Range("I3").Value = 0.04218439013093866
$ws.Range("J3").Value = 0.3550095850963846
$ws.Range("K3").Value = -0.5233929961551953

$ws.Range("D4").Value = -0.723599635428934
$ws.Range("E4").Value = -0.002672703882931809
$ws.Range("F4").Value = -0.7116246284967865
$ws.Range("G4").Value = -0.2195186361743094
$ws.Range("H4").Value = 0.2502861110986777
$ws.Range("I4").Value = 0.5631113060641236
$ws.Range("J4").Value = -0.3152912751874563
$ws.Range("K4").Value = -1.045203290365925

$ws.Range("C5").Value = -0.6071014650551335
$ws.Range("D5").Value = 0.1138254664908686
$ws.Range("E5").Value = -0.595126458122986
$ws.Range("F5").Value = -0.103020465800509
$ws.Range("G5").Value = 0.3667842814724781
$ws.Range("H5").Value = 0.6796094764379241
$ws.Range("I5").Value = -0.1987931048136559
$ws.Range("J5").Value = -0.9287051199921248
$ws.Range("K5").Value = 0.7334739763975026

$ws.Range("B6").Value = -0.5616080510579985
$ws.Range("C6").Value = 0.1593188804880037
$ws.Range("D6").Value = -0.549633044125851
$ws.Range("E6").Value = -0.05752705180337387
$ws.Range("F6").Value = 0.4122776954696132
$ws.Range("G6").Value = 0.7251028904350592
$ws.Range("H6").Value = -0.1532996908165208
$ws.Range("I6").Value = -0.8832117059949898
$ws.Range("J6").Value = 0.7789673903946376
$ws.Range("K6").Value = 0.1597481019993938

$ws.Range("B7").Value = 0.1417647591280393
$ws.Range("C7").Value = -0.5671871654858154
$ws.Range("D7").Value = -0.07508117316333833
$ws.Range("E7").Value = 0.3947235741096488
$ws.Range("F7").Value = 0.7075487690750948
$ws.Range("G7").Value = -0.1708538121764852
$ws.Range("H7").Value = -0.9007658273549541
$ws.Range("I7").Value = 0.7614132690346732
$ws.Range("J7").Value = 0.1421939806394294
$ws.Range("K7").Value = 0.379497744259143

$ws.Range("B8").Value = -0.4790798465348092
$ws.Range("C8").Value = 0.01302614578766791
$ws.Range("D8").Value = 0.482830893060655
$ws.Range("E8").Value = 0.795656088026101
$ws.Range("F8").Value = -0.082746493225479
$ws.Range("G8").Value = -0.8126585084039479
$ws.Range("H8").Value = 0.8495205879856794
$ws.Range("I8").Value = 0.2303012995904356
$ws.Range("J8").Value = 0.4676050632101492
$ws.Range("K8").Value = 0.3643276933347375

$ws.Range("B9").Value = 0.1916007792754515
$ws.Range("C9").Value = 0.6614055265484386
$ws.Range("D9").Value = 0.9742307215138846
$ws.Range("E9").Value = 0.0958281402623046
$ws.Range("F9").Value = -0.6340838749161644
$ws.Range("G9").Value = 1.028095221473463
$ws.Range("H9").Value = 0.4088759330782192
$ws.Range("I9").Value = 0.6461796966979327
$ws.Range("J9").Value = 0.542902326822521
$ws.Range("K9").Value = -0.2486961005069136

$ws.Range("B10").Value = 1.573432754301089
$ws.Range("C10").Value = 1.886257949266535
$ws.Range("D10").Value = 1.007855368014955
$ws.Range("E10").Value = 0.2779433528364856
$ws.Range("F10").Value = 1.940122449226113
$ws.Range("G10").Value = 1.320903160830869
$ws.Range("H10").Value = 1.558206924450583
$ws.Range("I10").Value = 1.454929554575171
$ws.Range("J10").Value = 0.6633311272457364
$ws.Range("K10").Value = 1.273820034913197

$ws.Range("B11").Value = 0.9422837133007778
$ws.Range("C11").Value = 0.06388113204919779
$ws.Range("D11").Value = -0.6660308831292712
$ws.Range("E11").Value = 0.9961482132603562
$ws.Range("F11").Value = 0.3769289248651124
$ws.Range("G11").Value = 0.6142326884848259
$ws.Range("H11").Value = 0.5109553186094142
$ws.Range("I11").Value = -0.2806431087200204
$ws.Range("J11").Value = 0.3298457989474406
$ws.Range("K11").Value = 0.1683237681281231

$ws.Range("B12").Value = 0.0678490295623069
$ws.Range("C12").Value = -0.6620629856161621
$ws.Range("D12").Value = 1.000116110773465
$ws.Range("E12").Value = 0.3808968223782215
$ws.Range("F12").Value = 0.6182005859979351
$ws.Range("G12").Value = 0.5149232161225235
$ws.Range("H12").Value = -0.2766752112069113
$ws.Range("I12").Value = 0.3338136964605497
$ws.Range("J12").Value = 0.1722916656412322

$ws.Range("B13").Value = -0.5264228954459207
$ws.Range("C13").Value = 1.135756200943707
$ws.Range("D13").Value = 0.5165369125484629
$ws.Range("E13").Value = 0.7538406761681764
$ws.Range("F13").Value = 0.6505633062927647
$ws.Range("G13").Value = -0.1410351210366699
$ws.Range("H13").Value = 0.4694537866307911
$ws.Range("I13").Value = 0.3079317558114735

$ws.Range("B14").Value = 0.894950019088042
$ws.Range("C14").Value = 0.2757307306927982
$ws.Range("D14").Value = 0.5130344943125118
$ws.Range("E14").Value = 0.4097571244371001
$ws.Range("F14").Value = -0.3818413028923346
$ws.Range("G14").Value = 0.2286476047751264
$ws.Range("H14").Value = 0.06712557395580883

$ws.Range("B15").Value = 0.2303995154407018
$ws.Range("C15").Value = 0.4677032790604154
$ws.Range("D15").Value = 0.3644259091850037
$ws.Range("E15").Value = -0.427172518144431
$ws.Range("F15").Value = 0.18331638952303
$ws.Range("G15").Value = 0.02179435870371246

$ws.Range("B16").Value = 0.4008418571243615
$ws.Range("C16").Value = 0.2975644872489498
$ws.Range("D16").Value = -0.4940339400804848
$ws.Range("E16").Value = 0.1164549675869761
$ws.Range("F16").Value = -0.04506706323234141

$ws.Range("B17").Value = 0.2679782848922332
$ws.Range("C17").Value = -0.5236201424372015
$ws.Range("D17").Value = 0.08686876523025952
$ws.Range("E17").Value = -0.07465326558905801

$ws.Range("B18").Value = -0.5417707991668423
$ws.Range("C18").Value = 0.06871810850061863
$ws.Range("D18").Value = -0.0928039223186989

$ws.Range("B19").Value = 0.0506862842519193
$ws.Range("C19").Value = -0.1108357465673982

$ws.Range("B20").Value = -0.1624199859130616
